# Apply the agenda-timing edits to the bepug meetup deck.
#
# Slide 3 ("Agenda part I") table: the "Welcome" row's start-time cell gets
# extra spacing inserted before the en-dash.
#
# Slide 4 ("Agenda part II") table: two previously-empty cells get their
# missing time ranges filled in, and the final row's time cell is updated
# from "20.15 - ..." to "21.00 - ....".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3 - Agenda part I
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tbl3 = $s3.Shapes.Item(2).Table

# Row 2, column 1 currently reads "…. – 18.30"; widen the gap before the dash.
$cell = $tbl3.Cell(2, 1)
$cell.Shape.TextFrame.TextRange.Text = "….     – 18.30"

# ---------------------------------------------------------------------
# Slide 4 - Agenda part II
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$tbl4 = $s4.Shapes.Item(2).Table

# Row 3, column 1 was empty -> "20.15 – 20.30"
$tbl4.Cell(3, 1).Shape.TextFrame.TextRange.Text = "20.15 – 20.30"

# Row 4, column 1 was empty -> "20.30 – 21.00"
$tbl4.Cell(4, 1).Shape.TextFrame.TextRange.Text = "20.30 – 21.00"

# Row 5, column 1 was "20.15 - …" -> "21.00 – …."
$tbl4.Cell(5, 1).Shape.TextFrame.TextRange.Text = "21.00 – …."
